$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 3, shifting old rows 3-5 down to 4-6.
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the new weekly record (same market/product metadata as surrounding rows).
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(3, 3).Value = "Los Lagos"
$ws.Cells.Item(3, 4).Value = 44495
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100107
$ws.Cells.Item(3, 8).Value = "Otros"
$ws.Cells.Item(3, 9).Value = 100107002
$ws.Cells.Item(3, 10).Value = "Chirimoya"
$ws.Cells.Item(3, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 270
$ws.Cells.Item(3, 14).Value = 19000
$ws.Cells.Item(3, 15).Value = 20000
$ws.Cells.Item(3, 16).Value = 19556
$ws.Cells.Item(3, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 2444
$ws.Cells.Item(3, 20).Value = 8
